$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.183.08'
$ws.Range('E2').Value = '  +0.85%  '

# Row 3
$ws.Range('D3').Value = '1.901.55'
$ws.Range('E3').Value = '  +0.88%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.08%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.81%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9995'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.05%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5212'
$ws.Range('D7').Style = 'Normal'

# Row 8
$ws.Range('E8').Value = '  +0.82%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07277'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.21%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.15'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.47%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9045'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.72%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08303'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +8.51%  '

# Row 13
$ws.Range('D13').Value = '1.914.08'
$ws.Range('E13').Value = '  +2.41%  '

# Row 14
$ws.Range('E14').Value = '  +3.47%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.289'
$ws.Range('D15').Style = 'Normal'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.09%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008666'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.24%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.18%  '

# Row 19
$ws.Range('E19').Value = '  +0.01%  '

# Row 20
$ws.Range('D20').Value = '27.222.88'
$ws.Range('E20').Value = '  +0.84%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.089'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.05%  '

# Row 22
$ws.Range('D22').Value = '2.165.58'
$ws.Range('E22').Value = '  +1.81%  '

# Row 23
$ws.Range('E23').Value = '  +0.91%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.438'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.90%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.320'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.33%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '146.38'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.15%  '

# Row 27
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.25'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.36%  '

# Row 28
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.747'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.19%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '115.01'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.05%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.838'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.43%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.894'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.13%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09269'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.04%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05083'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.11%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7982'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.42%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.245'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.35%  '

# Row 36
$ws.Range('E36').Value = '  +5.03%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.958'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.01%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.595'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.61%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5724'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.76%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02002'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.03%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.079'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.74%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.021'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.07%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.600'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.02%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '117.10'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.27%  '

# Row 45
$ws.Range('E45').Value = '  +1.40%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4862'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.23%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9996'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.04%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.08'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.81%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.631'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.07%  '

# Row 50
$ws.Range('E50').Value = '  +0.40%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.99'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.37%  '
